$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the date text (shared by A2:A4) ---
# The cells store a plain text date string ("2025-12-05" -> "2025-12-06"), not a real Excel date.
# Assigning a YYYY-MM-DD-looking string directly would get auto-converted to a date serial number,
# so we force Text format for the assignment and then clear the formatting again so the cells end
# up with no style applied (matching the original, unstyled data cells).
$ws.Range("A2:A4").NumberFormat = "@"
$ws.Range("A2").Value = "2025-12-06"
$ws.Range("A3").Value = "2025-12-06"
$ws.Range("A4").Value = "2025-12-06"
$ws.Range("A2:A4").ClearFormats()

# --- Row 2 (GLD) ---
$ws.Range("D2").Value = 390.28
$ws.Range("E2").Value = 73.59999999999999
$ws.Range("F2").Value = 0.62
$ws.Range("G2").Value = 50
$ws.Range("H2").Value = 70
$ws.Range("I2").Value = 80
$ws.Range("J2").Value = 93
$ws.Range("K2").Value = 62.5
$ws.Range("N2").Value = 51.54219175917372

# --- Row 3 (NEM) ---
$ws.Range("D3").Value = 91.31999999999999
$ws.Range("E3").Value = 58.6
$ws.Range("F3").Value = 0.66
$ws.Range("H3").Value = 76
$ws.Range("I3").Value = 70
$ws.Range("J3").Value = 83
$ws.Range("K3").Value = 61.5
$ws.Range("N3").Value = 51.54219175917372

# --- Row 4 (GC=F) ---
$ws.Range("D4").Value = 4270.1
$ws.Range("E4").Value = 73.8
$ws.Range("F4").Value = 1.23
$ws.Range("K4").Value = 59.7
$ws.Range("N4").Value = 51.54219175917372
